$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays text even for values that look numeric;
# Excel auto-converts numeric-looking strings assigned via .Value, which would
# change the cell type away from the original inline/shared text string type.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "67.817.63"
$ws.Range("E2").Value = "  -7.05%  "

# Row 3
$ws.Range("D3").Value = "3.698.24"
$ws.Range("E3").Value = "  -6.49%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "576.79"
$ws.Range("E5").Value = "  -5.65%  "

# Row 6
$ws.Range("D6").Value = "177.48"
$ws.Range("E6").Value = "  +7.17%  "

# Row 7
$ws.Range("D7").Value = "3.695.40"
$ws.Range("E7").Value = "  -6.37%  "

# Row 8
$ws.Range("D8").Value = "0.633"
$ws.Range("E8").Value = "  -6.69%  "

# Row 9
$ws.Range("D9").Value = "0.996"
$ws.Range("E9").Value = "  -0.39%  "

# Row 10
$ws.Range("D10").Value = "0.713"
$ws.Range("E10").Value = "  -5.23%  "

# Row 11
$ws.Range("D11").Value = "0.165"
$ws.Range("E11").Value = "  -9.77%  "

# Row 12
$ws.Range("D12").Value = "52.45"
$ws.Range("E12").Value = "  -5.98%  "

# Row 13
$ws.Range("D13").Value = "0.0000299"
$ws.Range("E13").Value = "  -9.65%  "

# Row 14
$ws.Range("D14").Value = "10.61"
$ws.Range("E14").Value = "  -4.32%  "

# Row 15
$ws.Range("D15").Value = "4.284.16"
$ws.Range("E15").Value = "  -6.56%  "

# Row 16
$ws.Range("D16").Value = "3.722.52"
$ws.Range("E16").Value = "  -5.93%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "19.39"
$ws.Range("E17").Value = "  -5.08%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.127"
$ws.Range("E18").Value = "  -3.23%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  -7.51%  "

# Row 20
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").Value = "1.14"
$ws.Range("E20").Value = "  -8.55%  "

# Row 21
$ws.Range("D21").Value = "67.750.50"
$ws.Range("E21").Value = "  -6.91%  "

# Row 22
$ws.Range("D22").Value = "408.32"
$ws.Range("E22").Value = "  -6.79%  "

# Row 23
$ws.Range("D23").Value = "4.58"
$ws.Range("E23").Value = "  -5.79%  "

# Row 24
$ws.Range("D24").Value = "88.38"
$ws.Range("E24").Value = "  -7.47%  "

# Row 25
$ws.Range("D25").Value = "3.09"
$ws.Range("E25").Value = "  -8.31%  "

# Row 26
$ws.Range("D26").Value = "12.82"
$ws.Range("E26").Value = "  -9.11%  "

# Row 27
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -2.84%  "

# Row 28
$ws.Range("E28").Value = "  -5.12%  "

# Row 29
$ws.Range("E29").Value = "  +0.19%  "

# Row 30
$ws.Range("D30").Value = "9.53"
$ws.Range("E30").Value = "  -8.72%  "

# Row 31
$ws.Range("D31").Value = "7.98"
$ws.Range("E31").Value = "  -0.44%  "

# Row 32
$ws.Range("D32").Value = "32.94"
$ws.Range("E32").Value = "  -8.29%  "

# Row 33
$ws.Range("D33").Value = "12.71"
$ws.Range("E33").Value = "  -6.45%  "

# Row 34
$ws.Range("E34").Value = "  -8.62%  "

# Row 35
$ws.Range("D35").Value = "44.35"

# Row 36
$ws.Range("D36").Value = "65.73"
$ws.Range("E36").Value = "  -6.38%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0917"
$ws.Range("E37").Value = "  -10.96%  "

# Row 38
$ws.Range("D38").Value = "600.45"
$ws.Range("E38").Value = "  -6.35%  "

# Row 39
$ws.Range("D39").Value = "0.401"
$ws.Range("E39").Value = "  -6.84%  "

# Row 40
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("D41").Value = "3.34"
$ws.Range("E41").Value = "  +16.30%  "

# Row 42
$ws.Range("E42").Value = "  -0.05%  "

# Row 44
$ws.Range("E44").Value = "  -11.14%  "

# Row 45
$ws.Range("E45").Value = "  -8.65%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.59"
$ws.Range("E46").Value = "  +0.09%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "9.46"
$ws.Range("E47").Value = "  -11.15%  "

# Row 48
$ws.Range("E48").Value = "  -8.84%  "

# Row 49
$ws.Range("D49").Value = "2.69"
$ws.Range("E49").Value = "  -15.72%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.737.46"
$ws.Range("E50").Value = "  -3.54%  "

# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "0.000268"
$ws.Range("E51").Value = "  -1.45%  "

# Restore the default cell style on column D so unaffected formatting metadata
# (e.g. no explicit style index) matches the original workbook.
$priceRange.Style = "Normal"
